$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NitroXHome")

# Fill in the 12th test case row (row 10) with the same values used for
# the 9th test case (Spot / Trader01@Tinyex / ETH / USDT)
$ws.Range("B10").Value = "Spot"
$ws.Range("C10").Value = "Trader01@Tinyex"
$ws.Range("D10").Value = "ETH"
$ws.Range("E10").Value = "USDT"

# Update selection to reflect where the user left the cursor after editing
$ws.Range("E10").Select()
